# 10.b.1.1 indicator workbook update
# - Update the header titles in row 1 (A1, C1) from "10.b.1" to "10.b.1.1"
# - Leave the selection on cell L8 (matches the saved view state in the target file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

$ws.Range("L8").Select()
